$wb = $excel.ActiveWorkbook

# --- Metadata sheet ---
$meta = $wb.Worksheets.Item("Metadata")

# "Date" property (A8/B8): refresh the generation timestamp
$meta.Range("B8").Value = "2025-06-13T15:45:04+00:00"

# "FHIR Version" property (A15/B15): target R4 instead of R4B
$meta.Range("B15").Value = "4.0.1"

# --- Elements sheet ---
$elem = $wb.Worksheets.Item("Elements")

$shortEle1 = "ele-1:All FHIR elements must have a @value or children {hasValue() or (children().count() > id.count())}" + [char]10 + "ext-1:Must have either extensions or value[x], not both {extension.exists() != value.exists()}"

# Row 2 = Extension (root) -> Constraint(s) column AJ: drop the
# "unless an empty Parameters resource ... or `$this is Parameters" clause
$elem.Range("AJ2").Value = $shortEle1

# Row 3 = Extension.id -> Type(s) column K: "id" becomes "string"
$elem.Range("K3").Value = "string" + [char]10

# Row 4 = Extension.extension -> Constraint(s) column AJ: same shortened wording
# (this cell used to duplicate the pre-edit long text that row 2 now matches)
$elem.Range("AJ4").Value = $shortEle1

# Row 6 = Extension.value[x] -> Definition column M: R4B link becomes R4
$elem.Range("M6").Value = "Value of extension - must be one of a constrained set of the data types (see [Extensibility](http://hl7.org/fhir/R4/extensibility.html) for a list)."

$wb.Save()
